$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 14399.667
$ws.Range("J51").Value = 19999.5
$ws.Range("L51").Value = 19999.5
$ws.Range("N51").Value = -20967.5
$ws.Range("H86").Value = 13592068
$ws.Range("I86").Value = 6950713.5
$ws.Range("J86").Value = 17861510
$ws.Range("K86").Value = 6950713.5
$ws.Range("L86").Value = 17861510
$ws.Range("M86").Value = -6949590.5
$ws.Range("N86").Value = -17863756
$ws.Range("H88").Value = 1978.238
$ws.Range("I88").Value = 1732.7
$ws.Range("J88").Value = 2201.4546
$ws.Range("K88").Value = 1732.7
$ws.Range("L88").Value = 2201.4546
$ws.Range("M88").Value = -1326.7
$ws.Range("N88").Value = -3013.4546
$ws.Range("H89").Value = 13592068
$ws.Range("I89").Value = 6950713.5
$ws.Range("J89").Value = 17861510
$ws.Range("K89").Value = 34753567.5
$ws.Range("L89").Value = 89307550
$ws.Range("M89").Value = -34747951.5
$ws.Range("N89").Value = -89318782
$ws.Range("H91").Value = 1978.238
$ws.Range("I91").Value = 1732.7
$ws.Range("J91").Value = 2201.4546
$ws.Range("K91").Value = 1732.7
$ws.Range("L91").Value = 2201.4546
$ws.Range("M91").Value = -328.7
$ws.Range("N91").Value = -5009.4546
$ws.Range("H98").Value = 1103.2424
$ws.Range("I98").Value = 1000.2414
$ws.Range("J98").Value = 1850
$ws.Range("K98").Value = 1000.2414
$ws.Range("L98").Value = 1850
$ws.Range("M98").Value = 497.7586
$ws.Range("N98").Value = -4846
$ws.Range("H122").Value = 1103.2424
$ws.Range("I122").Value = 1000.2414
$ws.Range("J122").Value = 1850
$ws.Range("K122").Value = 3000.7242
$ws.Range("L122").Value = 5550
$ws.Range("M122").Value = -550.7242000000001
$ws.Range("N122").Value = -10450
$ws.Range("H132").Value = 5968.34
$ws.Range("I132").Value = 5968.34
$ws.Range("K132").Value = 17905.02
$ws.Range("M132").Value = -15375.02
$ws.Range("H137").Value = 2392.7576
$ws.Range("I137").Value = 962.4706
$ws.Range("K137").Value = 2887.4118
$ws.Range("M137").Value = -337.4117999999999

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4216.6787
$ws.Range("I2").Value = 2179.087
$ws.Range("J2").Value = 13589.6
$ws.Range("K2").Value = 2179.087
$ws.Range("L2").Value = 13589.6
$ws.Range("M2").Value = -2066.087
$ws.Range("N2").Value = -13815.6
$ws.Range("H32").Value = 775.8219
$ws.Range("I32").Value = 703.2639
$ws.Range("K32").Value = 703.2639
$ws.Range("M32").Value = -416.2639
$ws.Range("H74").Value = 3365.9644
$ws.Range("I74").Value = 1431.3572
$ws.Range("K74").Value = 1431.3572
$ws.Range("M74").Value = -557.3571999999999
$ws.Range("H77").Value = 3365.9644
$ws.Range("I77").Value = 1431.3572
$ws.Range("K77").Value = 7156.786
$ws.Range("M77").Value = -2788.786
$ws.Range("H110").Value = 2019.0869
$ws.Range("I110").Value = 1232.1
$ws.Range("K110").Value = 1232.1
$ws.Range("M110").Value = 812.9000000000001
$ws.Range("H116").Value = 4216.6787
$ws.Range("I116").Value = 2179.087
$ws.Range("J116").Value = 13589.6
$ws.Range("K116").Value = 2179.087
$ws.Range("L116").Value = 13589.6
$ws.Range("M116").Value = 114.913
$ws.Range("N116").Value = -18177.6
$ws.Range("H132").Value = 3205.3928
$ws.Range("I132").Value = 3157.4443
$ws.Range("K132").Value = 9472.332900000001
$ws.Range("M132").Value = -6942.332900000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4216.6787
$ws.Range("I3").Value = 2179.087
$ws.Range("J3").Value = 13589.6
$ws.Range("K3").Value = 2179.087
$ws.Range("L3").Value = 13589.6
$ws.Range("M3").Value = -2065.087
$ws.Range("N3").Value = -13817.6
$ws.Range("H99").Value = 2905.5557
$ws.Range("I99").Value = 2358.3333
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 2358.3333
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -860.3332999999998
$ws.Range("N99").Value = -6996
$ws.Range("H105").Value = 2914.4092
$ws.Range("I105").Value = 3064.5789
$ws.Range("K105").Value = 3064.5789
$ws.Range("M105").Value = -1317.5789
$ws.Range("H107").Value = 1573.6818
$ws.Range("I107").Value = 1784.5555
$ws.Range("J107").Value = 1238.7646
$ws.Range("K107").Value = 1784.5555
$ws.Range("L107").Value = 1238.7646
$ws.Range("M107").Value = 135.4445000000001
$ws.Range("N107").Value = -5078.7646
$ws.Range("H134").Value = 1697
$ws.Range("I134").Value = 816.7
$ws.Range("J134").Value = 10500
$ws.Range("K134").Value = 2450.1
$ws.Range("L134").Value = 31500
$ws.Range("M134").Value = 84.89999999999964
$ws.Range("N134").Value = -36570

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3599.5
$ws.Range("I62").Value = 3600
$ws.Range("K62").Value = 3600
$ws.Range("M62").Value = -2976
$ws.Range("H65").Value = 3599.5
$ws.Range("I65").Value = 3600
$ws.Range("K65").Value = 18000
$ws.Range("M65").Value = -14880
$ws.Range("H105").Value = 42878.473
$ws.Range("I105").Value = 67137.664
$ws.Range("J105").Value = 1291.2858
$ws.Range("K105").Value = 67137.664
$ws.Range("L105").Value = 1291.2858
$ws.Range("M105").Value = -65390.664
$ws.Range("N105").Value = -4785.2858
$ws.Range("H107").Value = 723.25
$ws.Range("I107").Value = 649.75
$ws.Range("J107").Value = 870.25
$ws.Range("K107").Value = 649.75
$ws.Range("L107").Value = 870.25
$ws.Range("M107").Value = 1270.25
$ws.Range("N107").Value = -4710.25

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 284.58334
$ws.Range("J12").Value = 320
$ws.Range("L12").Value = 960
$ws.Range("N12").Value = -1306
$ws.Range("H17").Value = 333333630
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H55").Value = 6692.3
$ws.Range("J55").Value = 8646.429
$ws.Range("L55").Value = 25939.287
$ws.Range("N55").Value = -26293.287
$ws.Range("H68").Value = 13492.143
$ws.Range("J68").Value = 5444.654
$ws.Range("L68").Value = 16333.962
$ws.Range("N68").Value = -17955.962
$ws.Range("H71").Value = 13492.143
$ws.Range("J71").Value = 5444.654
$ws.Range("L71").Value = 49001.88600000001
$ws.Range("N71").Value = -57113.88600000001
$ws.Range("H107").Value = 504.66666
$ws.Range("J107").Value = 504.66666
$ws.Range("L107").Value = 1513.99998
$ws.Range("N107").Value = -5353.999980000001
$ws.Range("H113").Value = 448.72223
$ws.Range("J113").Value = 592.9
$ws.Range("L113").Value = 1778.7
$ws.Range("N113").Value = -6118.7

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17216.357
$ws.Range("J70").Value = 20803
$ws.Range("L70").Value = 20803
$ws.Range("N70").Value = -21343
$ws.Range("H73").Value = 17216.357
$ws.Range("J73").Value = 20803
$ws.Range("L73").Value = 20803
$ws.Range("N73").Value = -22675
$ws.Range("H132").Value = 3752.3438
$ws.Range("I132").Value = 3491.963
$ws.Range("J132").Value = 5158.4
$ws.Range("K132").Value = 10475.889
$ws.Range("L132").Value = 15475.2
$ws.Range("M132").Value = -7945.889000000001
$ws.Range("N132").Value = -20535.2

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 1847.8
$ws.Range("I35").Value = 1275.3334
$ws.Range("J35").Value = 7000
$ws.Range("K35").Value = 1275.3334
$ws.Range("L35").Value = 7000
$ws.Range("M35").Value = -939.3334
$ws.Range("N35").Value = -7672
$ws.Range("H46").Value = 4418.091
$ws.Range("I46").Value = 2499.6667
$ws.Range("J46").Value = 4721
$ws.Range("K46").Value = 2499.6667
$ws.Range("L46").Value = 4721
$ws.Range("M46").Value = -2311.6667
$ws.Range("N46").Value = -5097
$ws.Range("H82").Value = 1586.091
$ws.Range("J82").Value = 1570.6
$ws.Range("L82").Value = 1570.6
$ws.Range("N82").Value = -2292.6
$ws.Range("H85").Value = 1586.091
$ws.Range("J85").Value = 1570.6
$ws.Range("L85").Value = 1570.6
$ws.Range("N85").Value = -4066.6
$ws.Range("H122").Value = 77955.14
$ws.Range("I122").Value = 6318.778
$ws.Range("K122").Value = 18956.334
$ws.Range("M122").Value = -16506.334

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3602.6736
$ws.Range("I132").Value = 2347.625
$ws.Range("J132").Value = 9180.666999999999
$ws.Range("K132").Value = 7042.875
$ws.Range("L132").Value = 27542.001
$ws.Range("M132").Value = -4512.875
$ws.Range("N132").Value = -32602.001

Write-Host "Applied all changes"